$wb = $excel.ActiveWorkbook

# Update the survey sheet's "display.text" header to "display.prompt.text"
# to be compatible with rev 210 translations changes.
$surveySheet = $wb.Worksheets.Item("survey")
$surveySheet.Range("F1").Value = "display.prompt.text"

# Change config directory back to red cross demo: make the "survey" sheet
# the active sheet/tab again, with cell E7 selected.
$surveySheet.Activate()
$surveySheet.Range("E7").Select()
